$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1774
$ws.Range("F3").Value = 10395
$ws.Range("F7").Value = 83
$ws.Range("F8").Value = 1706
$ws.Range("F11").Value = 240
$ws.Range("F12").Value = 95
$ws.Range("F13").Value = 530
$ws.Range("F15").Value = 144
$ws.Range("F19").Value = 115
$ws.Range("F20").Value = 403
$ws.Range("F21").Value = 403
$ws.Range("F23").Value = 351
$ws.Range("F24").Value = 46
$ws.Range("F25").Value = 1049
$ws.Range("F26").Value = 1118
$ws.Range("F27").Value = 1202
$ws.Range("F29").Value = 1411
$ws.Range("F30").Value = 719
$ws.Range("F34").Value = 631
$ws.Range("F38").Value = 768
$ws.Range("F39").Value = 799
$ws.Range("F41").Value = 1274
$ws.Range("F42").Value = 832

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F9").Value = 8
$ws.Range("F11").Value = 71
$ws.Range("F14").Value = 116
$ws.Range("F19").Value = 1141
$ws.Range("G19").Value = 319
$ws.Range("F20").Value = 46
$ws.Range("F21").Value = 2242
$ws.Range("F22").Value = 1117
$ws.Range("F25").Value = 96
$ws.Range("F26").Value = 15
$ws.Range("E33").Value = "2024.07.12 19:30-07.12 22:00"
$ws.Range("F47").Value = 81

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 827
$ws.Range("F5").Value = 212
$ws.Range("F6").Value = 2569
$ws.Range("F7").Value = 4248
$ws.Range("F8").Value = 78
$ws.Range("F10").Value = 407
$ws.Range("F11").Value = 334
$ws.Range("F12").Value = 260
$ws.Range("F13").Value = 127
$ws.Range("F14").Value = 55

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1774
$ws.Range("F4").Value = 827
$ws.Range("F6").Value = 10395
$ws.Range("F7").Value = 212
$ws.Range("F8").Value = 4248
$ws.Range("F9").Value = 78
$ws.Range("F10").Value = 334
$ws.Range("F11").Value = 1706
$ws.Range("F13").Value = 240
$ws.Range("F14").Value = 8
$ws.Range("F16").Value = 144
$ws.Range("F17").Value = 116
$ws.Range("F20").Value = 403
$ws.Range("F21").Value = 403
$ws.Range("F22").Value = 46
$ws.Range("F23").Value = 351
$ws.Range("F24").Value = 46
$ws.Range("F25").Value = 2242
$ws.Range("F26").Value = 2242
$ws.Range("F27").Value = 1117
$ws.Range("F28").Value = 1049
$ws.Range("F29").Value = 1118
$ws.Range("F30").Value = 1202
$ws.Range("F31").Value = 96
$ws.Range("F32").Value = 1411
$ws.Range("F33").Value = 719
$ws.Range("F35").Value = 631
$ws.Range("F39").Value = 768
$ws.Range("E40").Value = "2024.07.12 19:30-07.12 22:00"
$ws.Range("F41").Value = 799
$ws.Range("F43").Value = 832
$ws.Range("F50").Value = 81
